# Generate Report for Handoff
# Update the "887e7144-0649-4393-a1aa-5989dfad35fd.md" row's timestamps
# on the Overview, zh-cn and de-de sheets to reflect a fresh handoff
# xliff generation / handoff datetime.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-09-02 04:44:49"

# --- zh-cn sheet: "Latest Handoff Datetime" column (H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-09-02 04:44:45"

# --- de-de sheet: "Latest Handoff Datetime" column (H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-09-02 04:44:49"
